$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-area/horizontal-scroll split ratio nudged slightly (983 -> 986 per mille)
$excel.ActiveWindow.TabRatio = 0.986

# Touch I1 and H3 so they materialize as blank (default-styled) cells,
# matching the empty neighbouring cells in their rows.
$ws.Range("I1").Font.Name = "Arial"
$ws.Range("I1").Font.Size = 10
$ws.Range("H3").Font.Name = "Arial"
$ws.Range("H3").Font.Size = 10

# New row 5 data: Network model version 1500327174.h5 (NVIDIA End-to-End Deep Learning)
$ws.Range("C5").Value = "1500327174.h5"
$ws.Range("D5").Value = 9.74
$ws.Range("E5").Value = 7.94
$ws.Range("F5").Value = 8.62
$ws.Range("G5").Value = 7.58
$ws.Range("H5").Value = 7.82
$ws.Range("I5").Formula = "=AVERAGE(D5:H5)"

# Update the view: selection to I11
$ws.Range("I11").Select()
